# Apply updated TPM-derived NATMI metrics to the LR-pairs sheet.
# Only numeric result columns (E..T) change; identifiers in A..D are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.311181333333333
$ws.Range("H2").Value = 6.933544
$ws.Range("I2").Value = 0.336172840858964
$ws.Range("J2").Value = 0.336172840858964
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.550723
$ws.Range("N2").Value = 34.652169
$ws.Range("O2").Value = 0.9771722872222799
$ws.Range("P2").Value = 0.9771722872222797
$ws.Range("Q2").Value = 26.695815384104
$ws.Range("R2").Value = 240.262338456936
$ws.Range("S2").Value = 0.3284987838041654
$ws.Range("T2").Value = 0.3284987838041653

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.311181333333333
$ws.Range("H3").Value = 6.933544
$ws.Range("I3").Value = 0.336172840858964
$ws.Range("J3").Value = 0.336172840858964
$ws.Range("O3").Value = 0.012666744083571
$ws.Range("P3").Value = 0.012666744083571
$ws.Range("Q3").Value = 0.3460485586773333
$ws.Range("R3").Value = 3.114437028096
$ws.Range("S3").Value = 0.004258215343007539
$ws.Range("T3").Value = 0.004258215343007538

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.311181333333333
$ws.Range("H4").Value = 6.933544
$ws.Range("I4").Value = 0.336172840858964
$ws.Range("J4").Value = 0.336172840858964
$ws.Range("M4").Value = 0.1201083333333333
$ws.Range("N4").Value = 0.360325
$ws.Range("O4").Value = 0.01016096869414922
$ws.Range("P4").Value = 0.01016096869414921
$ws.Range("Q4").Value = 0.2775921379777778
$ws.Range("R4").Value = 2.4983292418
$ws.Range("S4").Value = 0.00341584171179114
$ws.Range("T4").Value = 0.00341584171179114

# Row 5
$ws.Range("I5").Value = 0.5963918049111226
$ws.Range("J5").Value = 0.5963918049111226
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.550723
$ws.Range("N5").Value = 34.652169
$ws.Range("O5").Value = 0.9771722872222799
$ws.Range("P5").Value = 0.9771722872222797
$ws.Range("Q5").Value = 47.360058831104
$ws.Range("R5").Value = 426.240529479936
$ws.Range("S5").Value = 0.5827775440856254
$ws.Range("T5").Value = 0.5827775440856253

# Row 6
$ws.Range("I6").Value = 0.5963918049111226
$ws.Range("J6").Value = 0.5963918049111226
$ws.Range("O6").Value = 0.012666744083571
$ws.Range("P6").Value = 0.012666744083571
$ws.Range("S6").Value = 0.007554342366348196
$ws.Range("T6").Value = 0.007554342366348193

# Row 7
$ws.Range("I7").Value = 0.5963918049111226
$ws.Range("J7").Value = 0.5963918049111226
$ws.Range("M7").Value = 0.1201083333333333
$ws.Range("N7").Value = 0.360325
$ws.Range("O7").Value = 0.01016096869414922
$ws.Range("P7").Value = 0.01016096869414921
$ws.Range("Q7").Value = 0.4924659463111111
$ws.Range("R7").Value = 4.4321935168
$ws.Range("S7").Value = 0.006059918459149065
$ws.Range("T7").Value = 0.006059918459149063

# Row 8
$ws.Range("G8").Value = 0.4636166666666667
$ws.Range("H8").Value = 1.39085
$ws.Range("I8").Value = 0.06743535422991333
$ws.Range("J8").Value = 0.06743535422991333
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.550723
$ws.Range("N8").Value = 34.652169
$ws.Range("O8").Value = 0.9771722872222799
$ws.Range("P8").Value = 0.9771722872222797
$ws.Range("Q8").Value = 5.355107694850001
$ws.Range("R8").Value = 48.19596925365001
$ws.Range("S8").Value = 0.06589595933248905
$ws.Range("T8").Value = 0.06589595933248904

# Row 9
$ws.Range("G9").Value = 0.4636166666666667
$ws.Range("H9").Value = 1.39085
$ws.Range("I9").Value = 0.06743535422991333
$ws.Range("J9").Value = 0.06743535422991333
$ws.Range("O9").Value = 0.012666744083571
$ws.Range("P9").Value = 0.012666744083571
$ws.Range("Q9").Value = 0.06941639626666668
$ws.Range("R9").Value = 0.6247475664000001
$ws.Range("S9").Value = 0.0008541863742152696
$ws.Range("T9").Value = 0.0008541863742152694

# Row 10
$ws.Range("G10").Value = 0.4636166666666667
$ws.Range("H10").Value = 1.39085
$ws.Range("I10").Value = 0.06743535422991333
$ws.Range("J10").Value = 0.06743535422991333
$ws.Range("M10").Value = 0.1201083333333333
$ws.Range("N10").Value = 0.360325
$ws.Range("O10").Value = 0.01016096869414922
$ws.Range("P10").Value = 0.01016096869414921
$ws.Range("Q10").Value = 0.0556842251388889
$ws.Range("R10").Value = 0.5011580262500001
$ws.Range("S10").Value = 0.0006852085232090124
$ws.Range("T10").Value = 0.0006852085232090121
